$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 379.16666
$ws.Range("I12").Value = 379.16666
$ws.Range("K12").Value = 379.16666
$ws.Range("M12").Value = -209.16666
$ws.Range("H17").Value = 1990
$ws.Range("J17").Value = 1990
$ws.Range("L17").Value = 5970
$ws.Range("N17").Value = -6306
$ws.Range("H19").Value = 489.4
$ws.Range("I19").Value = 543
$ws.Range("K19").Value = 543
$ws.Range("M19").Value = -368
$ws.Range("H32").Value = 5886492.5
$ws.Range("I32").Value = 7566.6665
$ws.Range("J32").Value = 7146262
$ws.Range("K32").Value = 7566.6665
$ws.Range("L32").Value = 7146262
$ws.Range("M32").Value = -7240.6665
$ws.Range("N32").Value = -7146914
$ws.Range("H33").Value = 2916.8333
$ws.Range("H61").Value = 372.5
$ws.Range("I61").Value = 372.5
$ws.Range("K61").Value = 1117.5
$ws.Range("M61").Value = -945.5
$ws.Range("H69").Value = 26177.73
$ws.Range("I69").Value = 14332.667
$ws.Range("J69").Value = 27722.738
$ws.Range("K69").Value = 42998.001
$ws.Range("L69").Value = 83168.21400000001
$ws.Range("M69").Value = -42124.001
$ws.Range("N69").Value = -84916.21400000001
$ws.Range("H70").Value = 2133.9143
$ws.Range("I70").Value = 1623.8
$ws.Range("K70").Value = 4871.4
$ws.Range("M70").Value = -4601.4
$ws.Range("H72").Value = 26177.73
$ws.Range("I72").Value = 14332.667
$ws.Range("J72").Value = 27722.738
$ws.Range("K72").Value = 128994.003
$ws.Range("L72").Value = 249504.642
$ws.Range("M72").Value = -124626.003
$ws.Range("N72").Value = -258240.642
$ws.Range("H73").Value = 2133.9143
$ws.Range("I73").Value = 1623.8
$ws.Range("K73").Value = 4871.4
$ws.Range("M73").Value = -3935.4
$ws.Range("H80").Value = 2830.2693
$ws.Range("I80").Value = 588
$ws.Range("K80").Value = 1764
$ws.Range("M80").Value = -766
$ws.Range("H83").Value = 2830.2693
$ws.Range("I83").Value = 588
$ws.Range("K83").Value = 5292
$ws.Range("M83").Value = -300
$ws.Range("H88").Value = 1638.1428
$ws.Range("I88").Value = 4027.1667
$ws.Range("J88").Value = 682.5333000000001
$ws.Range("K88").Value = 4027.1667
$ws.Range("L88").Value = 682.5333000000001
$ws.Range("M88").Value = -3621.1667
$ws.Range("N88").Value = -1494.5333
$ws.Range("H91").Value = 1638.1428
$ws.Range("I91").Value = 4027.1667
$ws.Range("J91").Value = 682.5333000000001
$ws.Range("K91").Value = 4027.1667
$ws.Range("L91").Value = 682.5333000000001
$ws.Range("M91").Value = -2623.1667
$ws.Range("N91").Value = -3490.5333

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9895.291999999999
$ws.Range("I32").Value = 8564.522999999999
$ws.Range("K32").Value = 8564.522999999999
$ws.Range("M32").Value = -8277.522999999999
$ws.Range("H122").Value = 3535.5
$ws.Range("I122").Value = 1650.9678
$ws.Range("K122").Value = 4952.903399999999
$ws.Range("M122").Value = -2502.903399999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1810.25
$ws.Range("I20").Value = 1858.9412
$ws.Range("K20").Value = 1858.9412
$ws.Range("M20").Value = -1611.9412
$ws.Range("H86").Value = 2673.75
$ws.Range("I86").Value = 1522.3684
$ws.Range("J86").Value = 3960.5881
$ws.Range("K86").Value = 1522.3684
$ws.Range("L86").Value = 3960.5881
$ws.Range("M86").Value = -399.3684000000001
$ws.Range("N86").Value = -6206.5881
$ws.Range("H89").Value = 2673.75
$ws.Range("I89").Value = 1522.3684
$ws.Range("J89").Value = 3960.5881
$ws.Range("K89").Value = 7611.842000000001
$ws.Range("L89").Value = 19802.9405
$ws.Range("M89").Value = -1995.842000000001
$ws.Range("N89").Value = -31034.9405

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1427.75
$ws.Range("I16").Value = 903.6667
$ws.Range("K16").Value = 903.6667
$ws.Range("M16").Value = -616.6667
$ws.Range("H86").Value = 36750
$ws.Range("I86").Value = 48040.43
$ws.Range("J86").Value = 23577.834
$ws.Range("K86").Value = 48040.43
$ws.Range("L86").Value = 23577.834
$ws.Range("M86").Value = -46917.43
$ws.Range("N86").Value = -25823.834
$ws.Range("H89").Value = 36750
$ws.Range("I89").Value = 48040.43
$ws.Range("J89").Value = 23577.834
$ws.Range("K89").Value = 240202.15
$ws.Range("L89").Value = 117889.17
$ws.Range("M89").Value = -234586.15
$ws.Range("N89").Value = -129121.17
$ws.Range("H105").Value = 858.05884
$ws.Range("I105").Value = 830.4167
$ws.Range("J105").Value = 924.4
$ws.Range("K105").Value = 830.4167
$ws.Range("L105").Value = 924.4
$ws.Range("M105").Value = 916.5833
$ws.Range("N105").Value = -4418.4
$ws.Range("H113").Value = 1427.75
$ws.Range("I113").Value = 903.6667
$ws.Range("K113").Value = 903.6667
$ws.Range("M113").Value = 1266.3333
$ws.Range("H122").Value = 331133.3
$ws.Range("I122").Value = 394099.38
$ws.Range("K122").Value = 1182298.14
$ws.Range("M122").Value = -1179848.14

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 107.73077
$ws.Range("J2").Value = 208.83333
$ws.Range("L2").Value = 1252.99998
$ws.Range("N2").Value = -1478.99998
$ws.Range("H12").Value = 159.11539
$ws.Range("J12").Value = 183.94737
$ws.Range("L12").Value = 551.84211
$ws.Range("N12").Value = -897.84211
$ws.Range("H17").Value = 3367.3333
$ws.Range("I17").Value = 3367.3333
$ws.Range("K17").Value = 10101.9999
$ws.Range("M17").Value = -9932.999899999999
$ws.Range("H23").Value = 3965
$ws.Range("J23").Value = 947
$ws.Range("L23").Value = 2841
$ws.Range("N23").Value = -3311
$ws.Range("H34").Value = 754.5714
$ws.Range("I34").Value = 353
$ws.Range("J34").Value = 1055.75
$ws.Range("K34").Value = 1059
$ws.Range("L34").Value = 3167.25
$ws.Range("M34").Value = -975
$ws.Range("N34").Value = -3335.25
$ws.Range("H39").Value = 4039.7144
$ws.Range("I39").Value = 400
$ws.Range("J39").Value = 4319.6924
$ws.Range("K39").Value = 1200
$ws.Range("L39").Value = 12959.0772
$ws.Range("M39").Value = -906
$ws.Range("N39").Value = -13547.0772
$ws.Range("H50").Value = 1416.6666
$ws.Range("I50").Value = 150
$ws.Range("J50").Value = 2050
$ws.Range("K50").Value = 450
$ws.Range("L50").Value = 6150
$ws.Range("M50").Value = 31
$ws.Range("N50").Value = -7112
$ws.Range("H53").Value = 1416.6666
$ws.Range("I53").Value = 150
$ws.Range("J53").Value = 2050
$ws.Range("K53").Value = 450
$ws.Range("L53").Value = 6150
$ws.Range("M53").Value = 31
$ws.Range("N53").Value = -7112
$ws.Range("H55").Value = 3677.2856
$ws.Range("J55").Value = 4618
$ws.Range("L55").Value = 13854
$ws.Range("N55").Value = -14208
$ws.Range("H86").Value = 1660
$ws.Range("J86").Value = 1660
$ws.Range("L86").Value = 4980
$ws.Range("N86").Value = -7352
$ws.Range("H89").Value = 1660
$ws.Range("J89").Value = 1660
$ws.Range("L89").Value = 14940
$ws.Range("N89").Value = -26796
$ws.Range("H103").Value = 4998
$ws.Range("I103").Value = 4996.6665
$ws.Range("J103").Value = 5000
$ws.Range("K103").Value = 14989.9995
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = -14110.9995
$ws.Range("N103").Value = -16758

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 408.74075
$ws.Range("I97").Value = 329.94116
$ws.Range("J97").Value = 542.7
$ws.Range("K97").Value = 329.94116
$ws.Range("L97").Value = 542.7
$ws.Range("M97").Value = 166.05884
$ws.Range("N97").Value = -1534.7

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -705
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -893
$ws.Range("N27").ClearContents()
$ws.Range("H46").Value = 7779.0527
$ws.Range("I46").Value = 1672.2
$ws.Range("K46").Value = 1672.2
$ws.Range("M46").Value = -1484.2
$ws.Range("H68").Value = 4639.3057
$ws.Range("I68").Value = 3460.077
$ws.Range("J68").Value = 5305.826
$ws.Range("K68").Value = 3460.077
$ws.Range("L68").Value = 5305.826
$ws.Range("M68").Value = -2711.077
$ws.Range("N68").Value = -6803.826
$ws.Range("H71").Value = 4639.3057
$ws.Range("I71").Value = 3460.077
$ws.Range("J71").Value = 5305.826
$ws.Range("K71").Value = 17300.385
$ws.Range("L71").Value = 26529.13
$ws.Range("M71").Value = -13556.385
$ws.Range("N71").Value = -34017.13

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5272.643
$ws.Range("I81").Value = 6440.381
$ws.Range("J81").Value = 1769.4286
$ws.Range("K81").Value = 12880.762
$ws.Range("L81").Value = 3538.8572
$ws.Range("M81").Value = -11819.762
$ws.Range("N81").Value = -5660.8572
$ws.Range("H84").Value = 5272.643
$ws.Range("I84").Value = 6440.381
$ws.Range("J84").Value = 1769.4286
$ws.Range("K84").Value = 64403.81
$ws.Range("L84").Value = 17694.286
$ws.Range("M84").Value = -59099.81
$ws.Range("N84").Value = -28302.286
$ws.Range("H100").Value = 679.6667
$ws.Range("I100").Value = 673.65
$ws.Range("J100").Value = 800
$ws.Range("K100").Value = 1347.3
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = -806.3
$ws.Range("N100").Value = -2682
$ws.Range("H113").Value = 434.1
$ws.Range("I113").Value = 356.35294
$ws.Range("J113").Value = 874.6667
$ws.Range("K113").Value = 1069.05882
$ws.Range("L113").Value = 2624.0001
$ws.Range("M113").Value = 1100.94118
$ws.Range("N113").Value = -6964.0001
